$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 967:968, shifting the existing data (old rows 967-1091)
# down to new rows 969-1093. Excel COM Insert() copies formatting from the row
# above, which keeps the date number format on column D.
$ws.Rows("967:968").Insert()

# Row 967 - "Pintón" quality entry
$ws.Range("A967").Value2 = 7
$ws.Range("B967").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C967").Value2 = "Ñuble"
$ws.Range("D967").Value2 = 45212
$ws.Range("E967").Value2 = 16
$ws.Range("F967").Value2 = "Fruta"
$ws.Range("G967").Value2 = 100108
$ws.Range("H967").Value2 = "Tropicales y subtropicales"
$ws.Range("I967").Value2 = 100108006
$ws.Range("J967").Value2 = "Plátano"
$ws.Range("K967").Value2 = "Sin especificar"
$ws.Range("L967").Value2 = "Pintón"
$ws.Range("M967").Value2 = 180
$ws.Range("N967").Value2 = 23000
$ws.Range("O967").Value2 = 23000
$ws.Range("P967").Value2 = 23000
$ws.Range("Q967").Value2 = "$/caja 20 kilos"
$ws.Range("R967").Value2 = "Ecuador"
$ws.Range("S967").Value2 = 1150
$ws.Range("T967").Value2 = 20

# Row 968 - "Primera Pintón" quality entry
$ws.Range("A968").Value2 = 7
$ws.Range("B968").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C968").Value2 = "Ñuble"
$ws.Range("D968").Value2 = 45212
$ws.Range("E968").Value2 = 16
$ws.Range("F968").Value2 = "Fruta"
$ws.Range("G968").Value2 = 100108
$ws.Range("H968").Value2 = "Tropicales y subtropicales"
$ws.Range("I968").Value2 = 100108006
$ws.Range("J968").Value2 = "Plátano"
$ws.Range("K968").Value2 = "Sin especificar"
$ws.Range("L968").Value2 = "Primera Pintón"
$ws.Range("M968").Value2 = 150
$ws.Range("N968").Value2 = 24000
$ws.Range("O968").Value2 = 24000
$ws.Range("P968").Value2 = 24000
$ws.Range("Q968").Value2 = "$/caja 20 kilos"
$ws.Range("R968").Value2 = "Ecuador"
$ws.Range("S968").Value2 = 1200
$ws.Range("T968").Value2 = 20
